# Update price list values and date on the first sheet (Hoja1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the date in A1 (serial date value)
$ws.Range("A1").Value = 45436

# Update "Para CARGA" price column (D23:D27)
$ws.Range("D23").Value = 1931
$ws.Range("D24").Value = 2106
$ws.Range("D25").Value = 2250
$ws.Range("D26").Value = 2320
$ws.Range("D27").Value = 2790

# Update "Para DESCARGA" price column (D34:D35)
$ws.Range("D34").Value = 1520
$ws.Range("D35").Value = 1651
